$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I5").Value = 95.63254537748651
$ws.Range("I6").Value = 95.97187477484948
$ws.Range("I7").Value = 96.37304842850502

$ws.Range("G20").Value = 97.36876325527616
$ws.Range("G21").Value = 97.48742861060276
$ws.Range("G22").Value = 98.17268840828255

$ws.Range("H23").Value = 98.54592461718237
$ws.Range("H24").Value = 98.58871988640556
$ws.Range("H25").Value = 98.80991174912796

$ws.Range("I28").Value = 95.86406533252683
$ws.Range("I29").Value = 95.30154292309886

$ws.Range("G38").Value = 95.46878737064264
$ws.Range("G39").Value = 99.09465765755139

$ws.Range("H40").Value = 97.8679578120345
$ws.Range("H41").Value = 99.25694077661305

$ws.Range("I44").Value = 96.82041423455104
$ws.Range("I45").Value = 94.2249971903003

$ws.Range("G54").Value = 97.06109523406683
$ws.Range("G55").Value = 98.26278400375924

$ws.Range("H56").Value = 97.1206651797401
$ws.Range("H57").Value = 98.31643848757544
